$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows just above the current row 791, shifting the
# remaining rows (old 791-882) down to 793-884.
$ws.Rows("791:792").Insert()

# New row 791: Navel Late / Primera, fecha 2023-09-25 (45194)
$ws.Range("A791").Value = 4
$ws.Range("B791").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C791").Value = "Los Lagos"
$ws.Range("D791").Value = 45194
$ws.Range("E791").Value = 10
$ws.Range("F791").Value = "Fruta"
$ws.Range("G791").Value = 100102
$ws.Range("H791").Value = "Cítricos"
$ws.Range("I791").Value = 100102005
$ws.Range("J791").Value = "Naranja"
$ws.Range("K791").Value = "Navel Late"
$ws.Range("L791").Value = "Primera"
$ws.Range("M791").Value = 200
$ws.Range("N791").Value = 20000
$ws.Range("O791").Value = 20000
$ws.Range("P791").Value = 20000
$ws.Range("Q791").Value = "$/caja 15 kilos empedrada"
$ws.Range("R791").Value = "Región de O'Higgins"
$ws.Range("S791").Value = 1333
$ws.Range("T791").Value = 15

# New row 792: Navel Late / Segunda, fecha 2023-09-25 (45194)
$ws.Range("A792").Value = 4
$ws.Range("B792").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C792").Value = "Los Lagos"
$ws.Range("D792").Value = 45194
$ws.Range("E792").Value = 10
$ws.Range("F792").Value = "Fruta"
$ws.Range("G792").Value = 100102
$ws.Range("H792").Value = "Cítricos"
$ws.Range("I792").Value = 100102005
$ws.Range("J792").Value = "Naranja"
$ws.Range("K792").Value = "Navel Late"
$ws.Range("L792").Value = "Segunda"
$ws.Range("M792").Value = 200
$ws.Range("N792").Value = 16000
$ws.Range("O792").Value = 16000
$ws.Range("P792").Value = 16000
$ws.Range("Q792").Value = "$/caja 15 kilos empedrada"
$ws.Range("R792").Value = "Región de O'Higgins"
$ws.Range("S792").Value = 1067
$ws.Range("T792").Value = 15
